$d = $word.ActiveDocument

# --- 1. Remove justification (w:jc="both") from paragraphs 2 and 3 ---
$p2 = $d.Paragraphs.Item(2)
$p2.Format.Alignment = 0   # wdAlignParagraphLeft -> default, drops <w:jc/>

$p3 = $d.Paragraphs.Item(3)
$p3.Format.Alignment = 0   # wdAlignParagraphLeft -> default, drops <w:jc/>

# --- 2. Split paragraph 2's run after "K2 Cycle 1" and move the
#        "_GoBack" bookmark there. Bookmarks.Add relocates an existing
#        bookmark of the same name (removing it from its old spot in
#        paragraph 3) and, because the split happens exactly at the
#        insertion point with no character edits, both halves keep
#        their original run-level formatting/rsid. ---
$find = $d.Content
$find.Find.Execute("K2 Cycle 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $find.End
$insertRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $insertRange) | Out-Null

# --- 3. The trailing half ("to search for planets...") must become a
#        plain, unadorned run (matching the target XML, which mints a
#        fresh <w:r> there). Round-trip its text through a placeholder
#        so the engine rebuilds that run from scratch. ---
$secondPart = " to search for planets, brown dwarfs and stellar variability. "
$full = $d.Content.Text
$idx2 = $full.IndexOf($secondPart)
$rngB = $d.Range($idx2, $idx2 + $secondPart.Length)
$placeholderB = "Q_PLACEHOLDER_Q"
$rngB.Text = $placeholderB
$rngB2 = $d.Range($idx2, $idx2 + $placeholderB.Length)
$rngB2.Text = $secondPart

# --- 4. Re-merge the now-adjacent "exist" / "ence" runs in paragraph 3
#        into a single run by round-tripping the text through a
#        placeholder, forcing the engine to rebuild a unified run ---
$mergedText = "existence planets in short period orbits around White Dwarfs. Any transiting objects that are detected could be prime targets for follow up observations by the JWST (Loeb & "
$full = $d.Content.Text
$startIdx = $full.IndexOf("Such a survey (combined")
$existIdx = $full.IndexOf("exist", $startIdx)
$endIdx = $existIdx + $mergedText.Length

$mergeRange = $d.Range($existIdx, $endIdx)
$placeholder = "X_PLACEHOLDER_X"
$mergeRange.Text = $placeholder

$mergeRange2 = $d.Range($existIdx, $existIdx + $placeholder.Length)
$mergeRange2.Text = $mergedText
